$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("J2").Value = 7707
$ws.Range("K2").Value = 1294
$ws.Range("K3").Value = 1209
$ws.Range("C4").Value = 1846
$ws.Range("K4").Value = 262
$ws.Range("K5").Value = 74
$ws.Range("K6").Value = 1563
$ws.Range("C7").Value = 28390
$ws.Range("J7").Value = 29264
$ws.Range("K7").Value = 4402

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("K2").Value = 28
$ws.Range("K8").Value = 253
$ws.Range("K11").Value = 94
$ws.Range("K12").Value = 6
$ws.Range("K15").Value = 37
$ws.Range("K18").Value = 34
$ws.Range("K20").Value = 106
$ws.Range("K23").Value = 40
$ws.Range("K24").Value = 16
$ws.Range("K25").Value = 21
$ws.Range("K29").Value = 206
$ws.Range("K33").Value = 182
$ws.Range("K34").Value = 28
$ws.Range("K36").Value = 49
$ws.Range("K37").Value = 145
$ws.Range("K40").Value = 8
$ws.Range("K42").Value = 146
$ws.Range("K44").Value = 44
$ws.Range("K48").Value = 48
$ws.Range("K52").Value = 120
$ws.Range("K53").Value = 60
$ws.Range("C63").Value = 275
$ws.Range("J63").Value = 90
$ws.Range("K63").Value = 14
$ws.Range("K76").Value = 59
$ws.Range("K77").Value = 33
$ws.Range("K83").Value = 88
$ws.Range("K85").Value = 222
$ws.Range("K86").Value = 30
$ws.Range("K88").Value = 56
$ws.Range("K89").Value = 63
$ws.Range("K90").Value = 40
$ws.Range("K92").Value = 17
$ws.Range("K93").Value = 19
$ws.Range("K94").Value = 55
$ws.Range("K95").Value = 78
$ws.Range("K97").Value = 35
$ws.Range("K98").Value = 27
$ws.Range("C101").Value = 28390
$ws.Range("J101").Value = 29264
$ws.Range("K101").Value = 4402

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("K2").Value = 30
$ws.Range("K3").Value = 21
$ws.Range("K7").Value = 94

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("K2").Value = 8
$ws.Range("K7").Value = 63

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("K2").Value = 84
$ws.Range("K3").Value = 72
$ws.Range("K4").Value = 13
$ws.Range("K6").Value = 50
$ws.Range("K7").Value = 222

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("K3").Value = 24
$ws.Range("K7").Value = 120

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("K3").Value = 13
$ws.Range("K6").Value = 31
$ws.Range("K7").Value = 60

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("K2").Value = 75
$ws.Range("K6").Value = 85
$ws.Range("K7").Value = 253

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("K2").Value = 42
$ws.Range("K3").Value = 29
$ws.Range("K7").Value = 88

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("K3").Value = 70
$ws.Range("K7").Value = 182

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("K2").Value = 27
$ws.Range("K4").Value = 6
$ws.Range("K7").Value = 78

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("K2").Value = 30
$ws.Range("K3").Value = 46
$ws.Range("K6").Value = 55
$ws.Range("K7").Value = 145

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("K2").Value = 54
$ws.Range("K3").Value = 64
$ws.Range("K4").Value = 8
$ws.Range("K7").Value = 206

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("K4").Value = 6
$ws.Range("K6").Value = 21
$ws.Range("K7").Value = 48

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("K2").Value = 7
$ws.Range("K3").Value = 16
$ws.Range("K7").Value = 44

$ws = $wb.Worksheets.Item('River North')
$ws.Range("K2").Value = 13
$ws.Range("K6").Value = 30
$ws.Range("K7").Value = 59

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("K2").Value = 35
$ws.Range("K6").Value = 60
$ws.Range("K7").Value = 146

$ws = $wb.Worksheets.Item('Dunning')
$ws.Range("K2").Value = 5
$ws.Range("K7").Value = 16

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("K3").Value = 10
$ws.Range("K7").Value = 40

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("K2").Value = 31
$ws.Range("K3").Value = 31
$ws.Range("K7").Value = 106

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("K4").Value = 5
$ws.Range("K7").Value = 34

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("K6").Value = 9
$ws.Range("K7").Value = 49

$ws = $wb.Worksheets.Item('West Lawn')
$ws.Range("K6").Value = 10
$ws.Range("K7").Value = 19

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range("K2").Value = 8
$ws.Range("K7").Value = 28

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("K6").Value = 25
$ws.Range("K7").Value = 55

$ws = $wb.Worksheets.Item('East Side')
$ws.Range("K4").Value = 3
$ws.Range("K7").Value = 21

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("K5").Value = 2
$ws.Range("K7").Value = 37

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("K6").Value = 20
$ws.Range("K7").Value = 27

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("K2").Value = 8
$ws.Range("K7").Value = 28

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("K6").Value = 22
$ws.Range("K7").Value = 35

$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Range("K2").Value = 3
$ws.Range("K7").Value = 17

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("K2").Value = 12
$ws.Range("K7").Value = 56

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("K3").Value = 5
$ws.Range("K7").Value = 30

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("K6").Value = 10
$ws.Range("K7").Value = 40

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range("K6").Value = 5
$ws.Range("K7").Value = 33

$ws = $wb.Worksheets.Item('Hegewisch')
$ws.Range("K2").Value = 4
$ws.Range("K7").Value = 8

$ws = $wb.Worksheets.Item('Beverly')
$ws.Range("K2").Value = 3
$ws.Range("K7").Value = 6
